$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.08358602051541
$ws.Range("D2").Value = 1.088105204457594
$ws.Range("E2").Value = 1.086493625664122
$ws.Range("F2").Value = 1.097393813311776
$ws.Range("I2").Value = 1.056052127074691
$ws.Range("J2").Value = 1.088449995836141
$ws.Range("K2").Value = 1.090758750118277
$ws.Range("L2").Value = 1.089151326698996
$ws.Range("M2").Value = 1.100023665195198
$ws.Range("N2").Value = 1.089995719677188
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.085080318391272
$ws.Range("D3").Value = 1.089457720808337
$ws.Range("E3").Value = 1.087820347667749
$ws.Range("F3").Value = 1.098789047612937
$ws.Range("I3").Value = 1.056438783533779
$ws.Range("J3").Value = 1.089604303861765
$ws.Range("K3").Value = 1.091929677190461
$ws.Range("L3").Value = 1.090296218589511
$ws.Range("M3").Value = 1.101238934153897
$ws.Range("N3").Value = 1.091151666952609
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.086046091485267
$ws.Range("D4").Value = 1.090332011744496
$ws.Range("E4").Value = 1.088677990542423
$ws.Range("F4").Value = 1.099691054386227
$ws.Range("I4").Value = 1.056686722622062
$ws.Range("J4").Value = 1.090349606420225
$ws.Range("K4").Value = 1.092685919658971
$ws.Range("L4").Value = 1.091035652433775
$ws.Range("M4").Value = 1.1020239485487
$ws.Range("N4").Value = 1.09189802792619
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.086451836400025
$ws.Range("D5").Value = 1.090699358519524
$ws.Range("E5").Value = 1.089038348525766
$ws.Range("F5").Value = 1.100070070808682
$ws.Range("I5").Value = 1.056790418544731
$ws.Range("J5").Value = 1.090662550335407
$ws.Range("K5").Value = 1.09300350727194
$ws.Range("L5").Value = 1.091346182651342
$ws.Range("M5").Value = 1.102353651038725
$ws.Range("N5").Value = 1.092211416257626
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.086519947333439
$ws.Range("D6").Value = 1.090761025773949
$ws.Range("E6").Value = 1.08909884291661
$ws.Range("F6").Value = 1.100133698500853
$ws.Range("I6").Value = 1.056807798069853
$ws.Range("J6").Value = 1.090715072771002
$ws.Range("K6").Value = 1.093056812016826
$ws.Range("L6").Value = 1.091398302949919
$ws.Range("M6").Value = 1.102408991050141
$ws.Range("N6").Value = 1.092264013281105
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.086051514106516
$ws.Range("D7").Value = 1.09033692104998
$ws.Range("E7").Value = 1.088682806424007
$ws.Range("F7").Value = 1.099696119551637
$ws.Range("I7").Value = 1.056688110322495
$ws.Range("J7").Value = 1.090353789484604
$ws.Range("K7").Value = 1.092690164597618
$ws.Range("L7").Value = 1.091039803036135
$ws.Range("M7").Value = 1.102028355291564
$ws.Range("N7").Value = 1.091902216931
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.084091264780169
$ws.Range("D8").Value = 1.088562478393537
$ws.Range("E8").Value = 1.08694217336931
$ws.Range("F8").Value = 1.097865508609804
$ws.Range("I8").Value = 1.05618326713726
$ws.Range("J8").Value = 1.088840436182927
$ws.Range("K8").Value = 1.091154768358551
$ws.Range("L8").Value = 1.089538538313161
$ws.Range("M8").Value = 1.100434652507645
$ws.Range("N8").Value = 1.090386714494051
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.080628062452898
$ws.Range("D9").Value = 1.085428739408804
$ws.Range("E9").Value = 1.083868344134119
$ws.Range("F9").Value = 1.094633360550034
$ws.Range("I9").Value = 1.055276327579737
$ws.Range("J9").Value = 1.086161161237161
$ws.Range("K9").Value = 1.088438079220531
$ws.Range("L9").Value = 1.086882288711836
$ws.Range("M9").Value = 1.097615827732148
$ws.Range("N9").Value = 1.087703634670544
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.078312820823831
$ws.Range("D10").Value = 1.083334596848708
$ws.Range("E10").Value = 1.081814374523159
$ws.Range("F10").Value = 1.092473987633172
$ws.Range("I10").Value = 1.054659925106126
$ws.Range("J10").Value = 1.084366245265808
$ws.Range("K10").Value = 1.086619183107788
$ws.Range("L10").Value = 1.085103891802952
$ws.Range("M10").Value = 1.095729245734002
$ws.Range("N10").Value = 1.085906169712633
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.077308679838063
$ws.Range("D11").Value = 1.082426557398782
$ws.Range("E11").Value = 1.080923787777576
$ws.Range("F11").Value = 1.091537790629553
$ws.Range("I11").Value = 1.054390196385926
$ws.Range("J11").Value = 1.08358688973217
$ws.Range("K11").Value = 1.085829672785208
$ws.Range("L11").Value = 1.08433196868657
$ws.Range("M11").Value = 1.09491052165712
$ws.Range("N11").Value = 1.085125707404743
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.07693544517249
$ws.Range("D12").Value = 1.082089075353647
$ws.Range("E12").Value = 1.080592797271718
$ws.Range("F12").Value = 1.091189862674251
$ws.Range("I12").Value = 1.054289580735966
$ws.Range("J12").Value = 1.083297074347091
$ws.Range("K12").Value = 1.085536119880432
$ws.Range("L12").Value = 1.084044956231955
$ws.Range("M12").Value = 1.094606131601889
$ws.Range("N12").Value = 1.084835480448575
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.07701551680258
$ws.Range("D13").Value = 1.082161475356267
$ws.Range("E13").Value = 1.080663804414724
$ws.Range("F13").Value = 1.091264502770814
$ws.Range("I13").Value = 1.054311182467796
$ws.Range("J13").Value = 1.083359255668478
$ws.Range("K13").Value = 1.085599101354028
$ws.Range("L13").Value = 1.084106534402378
$ws.Range("M13").Value = 1.094671437042637
$ws.Range("N13").Value = 1.084897750074571
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.077277833329467
$ws.Range("D14").Value = 1.082398665057042
$ws.Range("E14").Value = 1.080896431847359
$ws.Range("F14").Value = 1.091509034559833
$ws.Range("I14").Value = 1.054381888176954
$ws.Range("J14").Value = 1.083562940227431
$ws.Range("K14").Value = 1.085805413633448
$ws.Range("L14").Value = 1.084308250010909
$ws.Range("M14").Value = 1.094885366415201
$ws.Range("N14").Value = 1.085101723888961
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.077439421576515
$ws.Range("D15").Value = 1.082544779333922
$ws.Range("E15").Value = 1.081039736335484
$ws.Range("F15").Value = 1.091659674266879
$ws.Range("I15").Value = 1.054425395725622
$ws.Range("J15").Value = 1.083688393365071
$ws.Range("K15").Value = 1.085932490359059
$ws.Range("L15").Value = 1.084432495611116
$ws.Range("M15").Value = 1.095017138154205
$ws.Range("N15").Value = 1.085227355184445
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.078379427438812
$ws.Range("D16").Value = 1.08339483332343
$ws.Range("E16").Value = 1.081873453947364
$ws.Range("F16").Value = 1.092536094669129
$ws.Range("I16").Value = 1.054677766447429
$ws.Range("J16").Value = 1.084417922840393
$ws.Range("K16").Value = 1.086671539446518
$ws.Range("L16").Value = 1.085155082019369
$ws.Range("M16").Value = 1.095783542880887
$ws.Range("N16").Value = 1.085957920675301
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07896862751755
$ws.Range("D17").Value = 1.083927707538785
$ws.Range("E17").Value = 1.082396096320455
$ws.Range("F17").Value = 1.09308553122342
$ws.Range("I17").Value = 1.054835314640334
$ws.Range("J17").Value = 1.08487495905132
$ws.Range("K17").Value = 1.087134608342535
$ws.Range("L17").Value = 1.085607838012633
$ws.Range("M17").Value = 1.096263796534962
$ws.Range("N17").Value = 1.08641560593007
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.079312141584747
$ws.Range("D18").Value = 1.084238402718114
$ws.Range("E18").Value = 1.082700828979925
$ws.Range("F18").Value = 1.093405895396204
$ws.Range("I18").Value = 1.054926937677999
$ws.Range("J18").Value = 1.085141333964621
$ws.Range("K18").Value = 1.087404524088522
$ws.Range("L18").Value = 1.085871743306649
$ws.Range("M18").Value = 1.096543745000813
$ws.Range("N18").Value = 1.086682359126302
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.079429244712991
$ws.Range("D19").Value = 1.084344321350503
$ws.Range("E19").Value = 1.08280471548146
$ws.Range("F19").Value = 1.093515112343838
$ws.Range("I19").Value = 1.054958132674391
$ws.Range("J19").Value = 1.085232126060908
$ws.Range("K19").Value = 1.087496527259108
$ws.Range("L19").Value = 1.08596169789334
$ws.Range("M19").Value = 1.096639170622686
$ws.Range("N19").Value = 1.08677328015778
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.078905428155178
$ws.Range("D20").Value = 1.08387054777365
$ws.Range("E20").Value = 1.082340033760756
$ws.Range("F20").Value = 1.093026593601628
$ws.Range("I20").Value = 1.05481843937797
$ws.Range("J20").Value = 1.084825944800234
$ws.Range("K20").Value = 1.087084944558312
$ws.Range("L20").Value = 1.085559280198428
$ws.Range("M20").Value = 1.096212288043942
$ws.Range("N20").Value = 1.086366522073125
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.077200594668682
$ws.Range("D21").Value = 1.082328824055799
$ws.Range("E21").Value = 1.080827934075218
$ws.Range("F21").Value = 1.091437031136956
$ws.Range("I21").Value = 1.054361078897529
$ws.Range("J21").Value = 1.083502969303522
$ws.Range("K21").Value = 1.085744667960752
$ws.Range("L21").Value = 1.084248857727286
$ws.Range("M21").Value = 1.094822377308826
$ws.Range("N21").Value = 1.085041667799461
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.076127236382542
$ws.Range("D22").Value = 1.081358345448068
$ws.Range("E22").Value = 1.079876132438243
$ws.Range("F22").Value = 1.090436549912817
$ws.Range("I22").Value = 1.054071050043576
$ws.Range("J22").Value = 1.082669261280588
$ws.Range("K22").Value = 1.084900281224713
$ws.Range("L22").Value = 1.083423286545586
$ws.Range("M22").Value = 1.093946865378845
$ws.Range("N22").Value = 1.084206775815506
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.076696385062482
$ws.Range("D23").Value = 1.081872924245546
$ws.Range("E23").Value = 1.080380805309303
$ws.Range("F23").Value = 1.090967026543519
$ws.Range("I23").Value = 1.054225034574209
$ws.Range("J23").Value = 1.083111407788795
$ws.Range("K23").Value = 1.085348069823939
$ws.Range("L23").Value = 1.083861096340235
$ws.Range("M23").Value = 1.094411146311561
$ws.Range("N23").Value = 1.08464955022247
$ws.Range("B24").Value = 1.019999999999999
$ws.Range("C24").Value = 1.078933985710144
$ws.Range("D24").Value = 1.083896376186362
$ws.Range("E24").Value = 1.082365366376437
$ws.Range("F24").Value = 1.093053225325105
$ws.Range("I24").Value = 1.054826065423736
$ws.Range("J24").Value = 1.084848092869106
$ws.Range("K24").Value = 1.087107386055062
$ws.Range("L24").Value = 1.085581221939818
$ws.Range("M24").Value = 1.096235563056059
$ws.Range("N24").Value = 1.086388701594796
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.081524489341544
$ws.Range("D25").Value = 1.086239741212611
$ws.Range("E25").Value = 1.084663817428997
$ws.Range("F25").Value = 1.095469736850398
$ws.Range("I25").Value = 1.055512860187475
$ws.Range("J25").Value = 1.086855335643543
$ws.Range("K25").Value = 1.088438079220531
$ws.Range("L25").Value = 1.087570305234925
$ws.Range("M25").Value = 1.098345838002558
$ws.Range("N25").Value = 1.088398794884206

Write-Host "Applied 264 cell updates"
